# Delete rows 2-7 (years 2004年-2009年), which shifts the former rows 8-9
# (2010年, 2011年) up to become the new rows 2-3.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2:G7").EntireRow.Delete()
